# VerkefnataflaBurnDown.xlsx edit script
# - Rename "Sheet3" -> "Fallaforritun" and move it between "Notendasögur" and "Burn-Down"
# - Fill the new "Fallaforritun" sheet with a small table describing the Loan classes
# - Turn that range into an Excel Table (Table1, TableStyleLight17)
# - Make "Fallaforritun" the active tab
# - Update view/selection state on the other sheets

$wb = $excel.ActiveWorkbook

# --- Reorder / rename sheets -------------------------------------------------
$fallaforritun = $wb.Worksheets.Item("Sheet3")
$fallaforritun.Name = "Fallaforritun"
$burnDown = $wb.Worksheets.Item("Burn-Down")
$fallaforritun.Move($burnDown)

# Re-fetch a fresh reference after the move/rename (stale COM refs stop writing).
$fb = $wb.Worksheets.Item("Fallaforritun")

# --- Populate the Fallaforritun sheet ----------------------------------------
$fb.Range("A1").Value = "Klasi"
$fb.Range("B1").Value = "Lýsing"
$fb.Range("C1").Value = "Já"

$fb.Range("A2").Value = "AccountType.py"
$fb.Range("B2").Value = "Klasi sem heldur utan um mismunandi reikninga sem eru lesnir úr skrá."

$fb.Range("A3").Value = "calcLoanFun.py"

# B4 is written before B3 so the shared-string table picks up the plain
# string ahead of the rich-text one (matches the original authoring order).
$fb.Range("B4").Value = "sem segir til um hvað á eftir að borga af láninu"

$b3 = $fb.Range("B3")
$b3.Value = "Inniheldur Loan sem sér um helstu upplýsingar hvers láns og calcLoan"
$b3.Characters(1, 10).Font.Name = "Calibri"
$b3.Characters(1, 10).Font.Size = 11
$b3.Characters(11, 5).Font.Name = "Calibri"
$b3.Characters(11, 5).Font.Size = 11
$b3.Characters(11, 5).Font.Bold = $true
$b3.Characters(16, 45).Font.Name = "Calibri"
$b3.Characters(16, 45).Font.Size = 11
$b3.Characters(61, 8).Font.Name = "Calibri"
$b3.Characters(61, 8).Font.Size = 11
$b3.Characters(61, 8).Font.Bold = $true

# --- Column widths (best-effort match of the original layout) ---------------
$fb.Columns.Item(1).ColumnWidth = 14.42
$fb.Columns.Item(2).ColumnWidth = 64.25
$fb.Columns.Item(3).ColumnWidth = 10.09

# --- Turn A1:C41 into an Excel table -----------------------------------------
$tbl = $fb.ListObjects.Add(1, $fb.Range("A1:C41"), 0, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = "TableStyleLight17"

# --- Selection on the Fallaforritun sheet + make it the active tab ----------
$fb.Range("B3").Select()
$fb.Activate()

# --- Notendasögur sheet: selection moves from D15 to T25:U25 -----------------
$notendasogur = $wb.Worksheets.Item("Notendasögur")
$notendasogur.Activate()
$notendasogur.Range("T25:U25").Select()

# --- Burn-Down sheet: selection narrows from A1:Q19 to F22 --------------------
$burnDown = $wb.Worksheets.Item("Burn-Down")
$burnDown.Activate()
$burnDown.Range("F22").Select()

# --- Leave Fallaforritun as the active / selected sheet, matching activeTab=1
$fb = $wb.Worksheets.Item("Fallaforritun")
$fb.Activate()
